$wb = $excel.ActiveWorkbook

# --- FlightFinder sheet updates ---
$flightFinder = $wb.Worksheets.Item("FlightFinder")

# Correct the "class" (E column) value for row 3
$flightFinder.Range("E3").Value = 9

# Add an explanatory comment on E1 ("class" header) matching the B1 comment style
$excel.UserName = "Lynda Ademola"
$comment = $flightFinder.Range("E1").AddComment("value between 0 - 9")

# Selection moves from H7 to I1:J6 on FlightFinder (and it is no longer the active tab)
[void]$flightFinder.Select()
[void]$flightFinder.Range("I1:J6").Select()

# --- Add the new "Flight" worksheet after FlightFinder ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$flightSheet = $wb.Worksheets.Add($null, $lastSheet)
$flightSheet.Name = "Flight"

$flightSheet.Range("A1").Value = "depart"
$flightSheet.Range("B1").Value = "return"

$departReturn = @(
    @(2, 1),
    @(1, 3),
    @(0, 0),
    @(1, 2),
    @(3, 1)
)
for ($i = 0; $i -lt $departReturn.Count; $i++) {
    $row = $i + 2
    $flightSheet.Cells.Item($row, 1).Value = $departReturn[$i][0]
    $flightSheet.Cells.Item($row, 2).Value = $departReturn[$i][1]
}

# New sheet becomes the active tab/selection
[void]$flightSheet.Select()
[void]$flightSheet.Range("F27").Select()
